$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.436.67'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.573.08'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '291.62'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').Value = '0.3729'
$ws.Range('E7').Value = '  -0.83%  '
$ws.Range('D8').Value = '49.99'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.3397'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.07561'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '1.143'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '21.31'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').Value = '6.017'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').Value = '6.957'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '1.572.78'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = '0.00001122'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '90.91'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').Value = '0.06760'
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '6.300'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '16.31'
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').Value = '22.426.89'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').Value = '2.695'
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').Value = '20.10'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').Value = '148.53'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').Value = '5.006'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('D30').Value = '125.66'
$ws.Range('D31').Value = '1.748.19'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').Value = '1.058'
$ws.Range('E32').Value = '  +7.62%  '
$ws.Range('D33').Value = '6.193'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '9.810'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = '0.08390'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('E37').Value = '  +2.02%  '
$ws.Range('D38').Value = '0.02488'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').Value = '0.2297'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').Value = '0.06522'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('D41').Value = '5.485'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').Value = '11.29'
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('D43').Value = '0.6231'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '14.06'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').Value = '0.5810'
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('D48').Value = '129.42'
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('D49').Value = '2.076'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').Value = '1.223'
$ws.Range('E50').Value = '  -5.84%  '
$ws.Range('D51').Value = '0.07331'
$ws.Range('E51').Value = '  +0.08%  '
